# 24 mayıs verileri eklendi.
# Adds the 2020-05-24 COVID-19 Turkey data row (row 74) to the "data" sheet
# and grows the Table3 structured table to cover the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New day's figures: date, test, case, death, recovered
$ws.Range("A74").Value = 43975
$ws.Range("B74").Value = 24589
$ws.Range("C74").Value = 1141
$ws.Range("D74").Value = 32
$ws.Range("E74").Value = 1092

# Grow the worksheet table (Table3) so the new row is part of it.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E74")) | Out-Null

# Match the saved selection state (active cell E73).
$ws.Range("E73").Select() | Out-Null
